# Natmi following Dr Hou advice
# Update the Rspo2-Lgr5 LR-pairs sheet: the "Target cluster" breakdown
# changes from 2 groups (FAPs/sCs) to 4 groups (ECs/FAPs/M1/sCs), with
# refreshed statistics for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Target cluster = ECs -------------------------------------
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rspo2"
$ws.Range("C2").Value = "Lgr5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.876175666666667
$ws.Range("H2").Value = 5.628527
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.220366
$ws.Range("N2").Value = 0.661098
$ws.Range("O2").Value = 0.07392971980931644
$ws.Range("P2").Value = 0.07748491934499904
$ws.Range("Q2").Value = 0.4134453269606666
$ws.Range("R2").Value = 3.721007942646
$ws.Range("S2").Value = 0.07392971980931644
$ws.Range("T2").Value = 0.07748491934499904

# --- Row 3: Target cluster = FAPs -------------------------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo2"
$ws.Range("C3").Value = "Lgr5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.876175666666667
$ws.Range("H3").Value = 5.628527
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.321162
$ws.Range("N3").Value = 6.963486
$ws.Range("O3").Value = 0.778717480428163
$ws.Range("P3").Value = 0.816165154137556
$ws.Range("Q3").Value = 4.354907662791333
$ws.Range("R3").Value = 39.194168965122
$ws.Range("S3").Value = 0.778717480428163
$ws.Range("T3").Value = 0.816165154137556

# --- Row 4 (new): Target cluster = M1 ----------------------------------
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo2"
$ws.Range("C4").Value = "Lgr5"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.876175666666667
$ws.Range("H4").Value = 5.628527
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02892933333333333
$ws.Range("N4").Value = 0.086788
$ws.Range("O4").Value = 0.009705387889255384
$ws.Range("P4").Value = 0.01017210939998877
$ws.Range("Q4").Value = 0.05427651125288889
$ws.Range("R4").Value = 0.488488601276
$ws.Range("S4").Value = 0.009705387889255384
$ws.Range("T4").Value = 0.01017210939998877

# --- Row 5 (new): Target cluster = sCs ----------------------------------
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo2"
$ws.Range("C5").Value = "Lgr5"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.876175666666667
$ws.Range("H5").Value = 5.628527
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4102925
$ws.Range("N5").Value = 0.820585
$ws.Range("O5").Value = 0.1376474118732652
$ws.Range("P5").Value = 0.09617781711745617
$ws.Range("Q5").Value = 0.7697808047158333
$ws.Range("R5").Value = 4.618684828295
$ws.Range("S5").Value = 0.1376474118732652
$ws.Range("T5").Value = 0.09617781711745617
